# Update 'por municipio' figures on the report sheet with freshly pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO OTORRINO
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 5
$ws.Range("G2").Value = 2
$ws.Range("J2").Value = 3
$ws.Range("L2").Value = 7

# Row 3 - PACOTE PRÉ-OPERATÓRIO PEDIÁTRICO CIRURGIA GERAL
$ws.Range("B3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("I3").Value = 3
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0

# Row 5 - ADENOIDECTOMIA PEDIÁTRICO
$ws.Range("E5").Value = 1
$ws.Range("L5").Value = 1

# Row 6 - AMIGDALECTOMIA- PEDIATRICO
$ws.Range("C6").Value = 3

# Row 7 - AMIGDALECTOMIA COM ADENOIDECTOMIA - PEDIATRICO
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 4
$ws.Range("G7").Value = 2
$ws.Range("J7").Value = 4
$ws.Range("L7").Value = 6

# Row 10 - HERNIOPLASTIA INGUINAL (BILATERAL) - PEDIATRICO
$ws.Range("B10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("I10").Value = 1
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0

# Row 11 - HERNIOPLASTIA UMBILICAL - PEDIATRICO
$ws.Range("B11").Value = 0
$ws.Range("F11").Value = 2
$ws.Range("I11").Value = 1
$ws.Range("M11").Value = 0

# Row 12 - ORQUIDOPEXIA BILATERAL - PEDIATRICO
$ws.Range("G12").Value = 0
$ws.Range("M12").Value = 1

# Row 14 - CORRECAO DE HIPOSPADIA (1º TEMPO) - PEDIATRICO
$ws.Range("B14").Value = 0
$ws.Range("G14").Value = 0

# Row 16 - POSTECTOMIA - PEDIATRICO
$ws.Range("B16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1

# Row 17 - TOTAL
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 10
$ws.Range("G17").Value = 6
$ws.Range("I17").Value = 5
$ws.Range("J17").Value = 7
$ws.Range("K17").Value = 7
$ws.Range("L17").Value = 16
$ws.Range("M17").Value = 2
